$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before D. This shifts D:K -> E:L (values + formats),
#    and the freshly-inserted D column starts out blank/General-formatted.
$ws.Columns("D").Insert()

# 2. Copy the number formats from column E (which now holds what used to be in D)
#    into the new column D, so D matches the existing D/E formatting (date format
#    on row 7, #,##0 on the data rows) instead of defaulting to General.
$ws.Columns("E").Copy()
$ws.Columns("D").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 3. Populate the refreshed financial data. For every row listed below we set
#    literal values across D:J (7 columns) to the newly reported figures. The
#    former "D:J" figures that used to live in those columns have already been
#    shifted one column to the right by the insert above (old J -> new K), and
#    that inherited K column is intentionally left untouched.
$newData = @{
    7   = @(43465,43100,42735,42369,42004,41639,41274)
    8   = @(11761800,11231100,11009000,11374700,11071800,11280500,11771900)
    9   = @(7420800,7133600,6995600,7611600,7241300,7552100,25903400)
    10  = @(4341000,4097500,4013400,3763200,3830500,3728400,-14131500)
    14  = @(-507100,-117800,-48200,-394900,66200,-68400,1991500)
    15  = @(473500,508300,565500,587900,584600,603600,3823700)
    17  = @(9635700,9819700,9740000,10093400,10315600,10665600,13250700)
    18  = @(2126200,1411500,1269000,1281300,756200,614900,-1478800)
    20  = @(-13500,4500,67300,61700,158200,81900,243500)
    21  = @(2587300,1922100,1949200,1932400,1654500,1309700,1704600)
    22  = @(49400,85300,124500,136900,166100,163800,190700)
    23  = @(2063300,1330700,1211700,1206100,748400,532900,-1426000)
    24  = @(383700,237900,224400,178400,173900,157100,-167200)
    26  = @(1679600,1092800,987400,1027700,574500,375900,-1258900)
    27  = @(1677400,1091700,986200,1027700,574500,375900,-1258900)
    32  = @(13500,-4500,-67300,-61700,-158200,-81900,-243500)
    33  = @(1677400,1091700,986200,1027700,574500,375900,-1258900)
    35  = @(1677400,1091700,986200,1027700,574500,375900,-1258900)
    38  = @(43465,43100,42735,42369,42004,41639,41274)
    41  = @(996300,803300,1113000,1404700,785400,883000,545300)
    43  = @(1941000,1890600,1808600,3684600,1768300,1815400,4104200)
    44  = @(1842300,1470900,1510200,1543900,1521400,1488900,3114600)
    45  = @(264800,235600,265900,294000,329900,426400,508300)
    46  = @(5044500,4400400,4697800,4413900,4404900,4613600,4350000)
    47  = @(2475100,2265300,2221500,2743300,3219000,3326700,4078400)
    48  = @(6878900,6598400,7170600,7442200,6966400,7032600)
    49  = @(595800,589000,612600,757300,637300,614900,1706500)
    52  = @(709100,807800,905400,1035600,699000,792100,2800500)
    54  = @(15703400,14661000,15608000,15924400,15926600,16379900,16765900)
    57  = @(1469800,1308200,1115300,1028900,958200,1592100,1757000)
    58  = @(24700,348900,564400,278300,409500,629400,1306000)
    59  = @(755100,755100,911100,1330700,874000,162700,780900)
    60  = @(2249600,2412300,2590700,2155300,2241700,2384200,2302300)
    61  = @(844900,980600,2126200,3174100,3431000,3981900,4218700)
    62  = @(1616800,1548300,1649300,1681900,1861400,1649300,2426900)
    66  = @(4716800,4945700,6369500,7013600,7536400,8022200,8401500)
    72  = @(6340400,5344000,4750500,4324100,3591500,3532000,5938700)
    76  = @(10986500,9715300,9238500,8910800,8390200,8357700,8364400)
    80  = @(43465,43100,42735,42369,42004,41639,41274)
    81  = @(1677400,1091700,986200,1027700,574500,375900,-1258900)
    83  = @(473500,504900,611500,587900,738300,611500,2932900)
    89  = @(1560700,1748100,1891700,1329600,1392400,824700,1166900)
    91  = @(-340000,-342200,-393800,-484700,-424100,-378100,-425200)
    94  = @(-291700,-249100,-294000,-488100,-277100,-333200,-80800)
    96  = @(-687800,-568800,-448800,-418500,-357900,-355700,-353400)
    100 = @(-1076000,-1799700,-1185900,-925600,-1229700,-148100,-1115300)
    101 = @(0,-7900,-1100,1100,16800,-5600,0)
    102 = @(193000,-308500,410600,-83000,-97600,337700,-29200)
}

$cols = @("D","E","F","G","H","I","J")

foreach ($rowNum in $newData.Keys) {
    $vals = $newData[$rowNum]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $colLetter = $cols[$i]
        $ws.Range($colLetter + $rowNum).Value = $vals[$i]
    }
}
